$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS: updated balance figures ---
$wsGlobal = $wb.Worksheets.Item("GLOBAL RESULTS")
$wsGlobal.Range("C3").Value = 17.509083171496854
$wsGlobal.Range("C5").Value = -0.7618846114845372
$wsGlobal.Range("C7").Value = 40.990928860131916
$wsGlobal.Range("C9").Value = -19.54421255285085
$wsGlobal.Range("C13").Value = 16.602082538926844
$wsGlobal.Range("C15").Value = -0.9563473932163205
$wsGlobal.Range("C17").Value = 17.724133755846978
$wsGlobal.Range("C19").Value = -24.532660780436228
$wsGlobal.Range("C23").Value = 16.602082538926844
$wsGlobal.Range("C25").Value = -0.9563473932163205
$wsGlobal.Range("C27").Value = 17.724133755846978
$wsGlobal.Range("C29").Value = -24.532660780436228
$wsGlobal.Range("C33").Value = 16.602082538926844
$wsGlobal.Range("C35").Value = -0.9563473932163205
$wsGlobal.Range("C37").Value = 17.724133755846978
$wsGlobal.Range("C39").Value = -24.532660780436228
$wsGlobal.Range("C43").Value = 17.04956453589054
$wsGlobal.Range("C45").Value = -0.6732119944890558
$wsGlobal.Range("C47").Value = 29.20314662787694
$wsGlobal.Range("C49").Value = -17.269542021311445
$wsGlobal.Range("C53").Value = 16.956380611370676
$wsGlobal.Range("C55").Value = -0.7799135731491034
$wsGlobal.Range("C57").Value = 26.812749972172895
$wsGlobal.Range("C59").Value = -20.00669972422566
$wsGlobal.Range("C62").Value = 11.997784147183046
$wsGlobal.Range("C63").Value = 26.812749972172895
$wsGlobal.Range("C64").Value = 32.8363585964546
$wsGlobal.Range("C69").Value = 60154.70017136331
$wsGlobal.Range("C70").Value = 3037072.693850185
$wsGlobal.Range("C71").Value = 2976917.9936788203
$wsGlobal.Range("C76").Value = -6986.489026558615

# --- LANDING GEARS: updated balance figures ---
$wsLandingGears = $wb.Worksheets.Item("LANDING GEARS")
$wsLandingGears.Range("C5").Value = 13.534850162764503
$wsLandingGears.Range("C6").Value = 13.534850162764458
$wsLandingGears.Range("C7").Value = 16.43877328847313
$wsLandingGears.Range("C8").Value = 16.43877328847313
$wsLandingGears.Range("C9").Value = 16.43877328847313
$wsLandingGears.Range("C10").Value = 16.43877328847313
$wsLandingGears.Range("C23").Value = 16.43877328847313

# --- remove the SYSTEMS sheet entirely ---
$wsSystems = $wb.Worksheets.Item("SYSTEMS")
$wsSystems.Delete()

# keep the original active sheet selection (GLOBAL RESULTS)
$wb.Worksheets.Item("GLOBAL RESULTS").Activate()
